$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update goal values in column B
$ws.Range("B2").Value = 132
$ws.Range("B3").Value = -314
$ws.Range("B5").Value = 0.34
$ws.Range("B6").Value = 0.62
$ws.Range("B7").Value = 0.62
$ws.Range("B8").Value = 0.34

# Move active cell selection to B3
$ws.Range("B3").Select()
